$d = $word.ActiveDocument

# Step 1: remove the trailing "dev分支。" run's text. At this point in the
# document it is still the only occurrence of this phrase, so Find locates
# it unambiguously; deleting it first (before any insertion) avoids any
# risk of the following step's Find matching inside newly-inserted text.
$r1 = $d.Content
$r1.Find.Execute("dev分支。", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$r1.Delete()

# Step 2: the first run now ends the paragraph (right before the
# "_GoBack" bookmark). Collapse the found range to its end and insert the
# merged/extended text right after it, keeping everything in that same
# run and leaving the bookmark bracketing untouched.
$r2 = $d.Content
$r2.Find.Execute("创建了一个", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$r2.Collapse(0)
$r2.InsertAfter("dev分支。使用git创建分支简单又便捷。")
